# Add 2022-Q4 data:
#  1. "总计" (summary) sheet gets a new row inserted for 2022-Q4, pushing
#     2022-Q3 / 2022-Q2 / 2022-Q1 down by one row.
#  2. A brand-new "2022-Q4" worksheet is created (positioned right after
#     "总计", i.e. before "2022-Q3") holding the quarterly fund-holding table.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) "总计" sheet: shift existing data rows down and insert 2022-Q4 on
#    top. All values are written as literals (taken straight from the
#    known final state) so no floating point round-trip noise is
#    introduced.
# ---------------------------------------------------------------------
$summary = $wb.Worksheets.Item("总计")

# Give the about-to-be-used A5 cell the same style as the existing
# A2:A4 cells (centered/bold/bordered) before the value is (re)written.
$summary.Range("A2").Copy($summary.Range("A5"))

$summary.Cells.Item(2, 1).Value = 0
$summary.Cells.Item(2, 2).Value = "2022-Q4"
$summary.Cells.Item(2, 3).Value = 2
$summary.Cells.Item(2, 4).Value = 0.19

$summary.Cells.Item(3, 1).Value = 1
$summary.Cells.Item(3, 2).Value = "2022-Q3"
$summary.Cells.Item(3, 3).Value = 4
$summary.Cells.Item(3, 4).Value = 0.68

$summary.Cells.Item(4, 1).Value = 2
$summary.Cells.Item(4, 2).Value = "2022-Q2"
$summary.Cells.Item(4, 3).Value = 13
$summary.Cells.Item(4, 4).Value = 2.46

$summary.Cells.Item(5, 1).Value = 3
$summary.Cells.Item(5, 2).Value = "2022-Q1"
$summary.Cells.Item(5, 3).Value = 10
$summary.Cells.Item(5, 4).Value = 2.7

# ---------------------------------------------------------------------
# 2) New "2022-Q4" worksheet, inserted between "总计" and "2022-Q3".
# ---------------------------------------------------------------------
$q3 = $wb.Worksheets.Item("2022-Q3")
$q4 = $wb.Worksheets.Add($q3)
$q4.Name = "2022-Q4"

$q4.Range("B1").Value = "基金代码"
$q4.Range("C1").Value = "基金名称"
$q4.Range("D1").Value = "基金规模"
$q4.Range("E1").Value = "股票总仓位"
$q4.Range("F1").Value = "仓位占比"
$q4.Range("G1").Value = "持有市值(亿元)"
$q4.Range("H1").Value = "仓位排名"

# Match header styling (bold/centered/bordered) used on every other
# sheet's header row, then on the A-column index cells.
$q3.Range("B1:H1").Copy($q4.Range("B1:H1"))
$q3.Range("A2").Copy($q4.Range("A2"))
$q3.Range("A2").Copy($q4.Range("A3"))

$q4.Range("B1").Value = "基金代码"
$q4.Range("C1").Value = "基金名称"
$q4.Range("D1").Value = "基金规模"
$q4.Range("E1").Value = "股票总仓位"
$q4.Range("F1").Value = "仓位占比"
$q4.Range("G1").Value = "持有市值(亿元)"
$q4.Range("H1").Value = "仓位排名"

$q4.Cells.Item(2, 1).Value = 0
$q4.Cells.Item(2, 2).Value = "'012421"
$q4.Cells.Item(2, 3).Value = "华夏优加生活混合A"
$q4.Cells.Item(2, 4).Value = "'8.01"
$q4.Cells.Item(2, 5).Value = "'88.08"
$q4.Cells.Item(2, 6).Value = "'2.28"
$q4.Cells.Item(2, 7).Value = "'0.1826"
$q4.Cells.Item(2, 8).Value = 10

$q4.Cells.Item(3, 1).Value = 1
$q4.Cells.Item(3, 2).Value = "'012422"
$q4.Cells.Item(3, 3).Value = "华夏优加生活混合C"
$q4.Cells.Item(3, 4).Value = "'0.32"
$q4.Cells.Item(3, 5).Value = "'88.08"
$q4.Cells.Item(3, 6).Value = "'2.28"
$q4.Cells.Item(3, 7).Value = "'0.0073"
$q4.Cells.Item(3, 8).Value = 10
